# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns on the per-locale sheets now that the handback round-tripped back in sync with
# en-US, flips the Status everywhere it is shown, and widens the columns that now hold the
# longer text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$targetFileName = "6a1f3617-07d8-4b4c-add0-535573617591.md"
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4d91e845beee05203e7c143da2fb0dab45c7238a/e2e/6a1f3617-07d8-4b4c-add0-535573617591.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: flip the per-locale Status columns (E = zh-cn, F = de-de) ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Status column got wider to fit the longer text
$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 29.17

# Latest Target File (I) - new hyperlink to the handed-off markdown file
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $targetUrl, "", "", $targetFileName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $targetUrl, "", "", $targetFileName)
$wsZhCn.Columns.Item(9).ColumnWidth = 39.17
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17

# Latest Handback File (J) mirrors the Latest Handoff File (G) for that row
$wsZhCn.Range("J2").Value = $wsZhCn.Range("G2").Value2
$wsZhCn.Range("J3").Value = $wsZhCn.Range("G3").Value2

# Latest Handback DateTime (K)
$wsZhCn.Range("K2").Value = "2016-08-29 07:04:30"
$wsZhCn.Range("K3").Value = "2016-08-29 07:04:30"

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 29.17

# Latest Target File (I) - new hyperlink to the handed-off markdown file
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $targetUrl, "", "", $targetFileName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $targetUrl, "", "", $targetFileName)
$wsDeDe.Columns.Item(9).ColumnWidth = 39.17
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17

# Latest Handback File (J) mirrors the Latest Handoff File (G) for that row
$wsDeDe.Range("J2").Value = $wsDeDe.Range("G2").Value2
$wsDeDe.Range("J3").Value = $wsDeDe.Range("G3").Value2

# Latest Handback DateTime (K)
$wsDeDe.Range("K2").Value = "2016-08-29 07:04:37"
$wsDeDe.Range("K3").Value = "2016-08-29 07:04:37"
